$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Relabel the event form fields to lowercase (for the database)
$ws.Range("E11").Value = "host"
$ws.Range("E5").Value = "title"
$ws.Range("E6").Value = "date"
$ws.Range("E7").Value = "time"
$ws.Range("E8").Value = "skill level"
$ws.Range("E9").Value = "message"
$ws.Range("E10").Value = "admin"

# Add note about keeping first letter lowercase for the database
$ws.Range("I3").Value = "Keep first letter lowercase for database"

# Update the selected cell to reflect the newly added note
$ws.Range("I3").Select()
